$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1262
$ws.Range("J19").Value = 1243.7333
$ws.Range("L19").Value = 1243.7333
$ws.Range("N19").Value = -1593.7333
$ws.Range("H33").Value = 686.25
$ws.Range("I33").Value = 407.33334
$ws.Range("K33").Value = 407.33334
$ws.Range("M33").Value = -178.33334
$ws.Range("H98").Value = 1119.4412
$ws.Range("I98").Value = 1009.0345
$ws.Range("K98").Value = 1009.0345
$ws.Range("M98").Value = 488.9655
$ws.Range("H122").Value = 1119.4412
$ws.Range("I122").Value = 1009.0345
$ws.Range("K122").Value = 3027.1035
$ws.Range("M122").Value = -577.1035000000002
$ws.Range("H132").Value = 6088.306
$ws.Range("I132").Value = 6088.306
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18264.918
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15734.918
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 2453.75
$ws.Range("I137").Value = 994.86664
$ws.Range("K137").Value = 2984.59992
$ws.Range("M137").Value = -434.5999199999997

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 3323.2
$ws.Range("I22").Value = 1572
$ws.Range("J22").Value = 5950
$ws.Range("K22").Value = 1572
$ws.Range("L22").Value = 5950
$ws.Range("M22").Value = -1273
$ws.Range("N22").Value = -6548
$ws.Range("H41").Value = 1982.1428
$ws.Range("I41").Value = 1982.1428
$ws.Range("K41").Value = 1982.1428
$ws.Range("M41").Value = -1568.1428
$ws.Range("H45").Value = 3844
$ws.Range("I45").Value = 3962.5715
$ws.Range("K45").Value = 3962.5715
$ws.Range("M45").Value = -3585.5715
$ws.Range("H132").Value = 2807.0881
$ws.Range("I132").Value = 2755.7878
$ws.Range("K132").Value = 8267.3634
$ws.Range("M132").Value = -5737.3634

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 181.6
$ws.Range("I22").Value = 181.6
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 181.6
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -8.599999999999994
$ws.Range("N22").ClearContents()
$ws.Range("H25").Value = 157.71428
$ws.Range("J25").Value = 100
$ws.Range("L25").Value = 100
$ws.Range("N25").Value = -570
$ws.Range("H86").Value = 3408.9092
$ws.Range("I86").Value = 3655.3333
$ws.Range("K86").Value = 3655.3333
$ws.Range("M86").Value = -2532.3333
$ws.Range("H89").Value = 3408.9092
$ws.Range("I89").Value = 3655.3333
$ws.Range("K89").Value = 18276.6665
$ws.Range("M89").Value = -12660.6665
$ws.Range("H105").Value = 3095.85
$ws.Range("I105").Value = 3295.7058
$ws.Range("K105").Value = 3295.7058
$ws.Range("M105").Value = -1548.7058
$ws.Range("H134").Value = 1927.5555
$ws.Range("I134").Value = 864.4706
$ws.Range("K134").Value = 2593.4118
$ws.Range("M134").Value = -58.41179999999986

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 801.6070999999999
$ws.Range("I22").Value = 764.4167
$ws.Range("J22").Value = 1024.75
$ws.Range("K22").Value = 764.4167
$ws.Range("L22").Value = 1024.75
$ws.Range("M22").Value = -414.4167
$ws.Range("N22").Value = -1724.75
$ws.Range("H31").Value = 16399.8
$ws.Range("I31").Value = 1351.6471
$ws.Range("J31").Value = 36078.152
$ws.Range("K31").Value = 1351.6471
$ws.Range("L31").Value = 36078.152
$ws.Range("M31").Value = -1056.6471
$ws.Range("N31").Value = -36668.152
$ws.Range("H34").Value = 16399.8
$ws.Range("I34").Value = 1351.6471
$ws.Range("J34").Value = 36078.152
$ws.Range("K34").Value = 1351.6471
$ws.Range("L34").Value = 36078.152
$ws.Range("M34").Value = -1149.6471
$ws.Range("N34").Value = -36482.152
$ws.Range("H94").Value = 1989.2667
$ws.Range("I94").Value = 1509.75
$ws.Range("J94").Value = 2163.6365
$ws.Range("K94").Value = 1509.75
$ws.Range("L94").Value = 2163.6365
$ws.Range("M94").Value = -1058.75
$ws.Range("N94").Value = -3065.6365
$ws.Range("H132").Value = 3786.8333
$ws.Range("I132").Value = 3857.2727
$ws.Range("K132").Value = 11571.8181
$ws.Range("M132").Value = -9041.8181
$ws.Range("H134").Value = 1358.3846
$ws.Range("I134").Value = 1239
$ws.Range("K134").Value = 3717
$ws.Range("M134").Value = -1182

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 219.71428
$ws.Range("I33").Value = 136.33333
$ws.Range("J33").Value = 369.8
$ws.Range("K33").Value = 817.9999799999999
$ws.Range("L33").Value = 2218.8
$ws.Range("M33").Value = -534.9999799999999
$ws.Range("N33").Value = -2784.8
$ws.Range("H107").Value = 549.0714
$ws.Range("J107").Value = 549.0714
$ws.Range("L107").Value = 1647.2142
$ws.Range("N107").Value = -5487.2142
$ws.Range("H122").Value = 15390065
$ws.Range("I122").Value = 365.77777
$ws.Range("K122").Value = 3291.99993
$ws.Range("M122").Value = -841.9999299999999
$ws.Range("H134").Value = 1452.1765
$ws.Range("I134").Value = 691.9286
$ws.Range("K134").Value = 2075.7858
$ws.Range("M134").Value = 2994.2142
$ws.Range("H139").Value = 2285.1333
$ws.Range("I139").Value = 1939.75
$ws.Range("J139").Value = 3666.6667
$ws.Range("K139").Value = 5819.25
$ws.Range("L139").Value = 11000.0001
$ws.Range("M139").Value = -679.25
$ws.Range("N139").Value = -21280.0001
$ws.Range("H141").Value = 2777.9285
$ws.Range("I141").Value = 2068.5386
$ws.Range("J141").Value = 12000
$ws.Range("K141").Value = 6205.6158
$ws.Range("L141").Value = 36000
$ws.Range("M141").Value = -1025.6158
$ws.Range("N141").Value = -46360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18617.154
$ws.Range("I70").Value = 8249.75
$ws.Range("K70").Value = 8249.75
$ws.Range("M70").Value = -7979.75
$ws.Range("H73").Value = 18617.154
$ws.Range("I73").Value = 8249.75
$ws.Range("K73").Value = 8249.75
$ws.Range("M73").Value = -7313.75
$ws.Range("H126").Value = 3882.4736
$ws.Range("I126").Value = 3247.25
$ws.Range("J126").Value = 4971.4287
$ws.Range("K126").Value = 9741.75
$ws.Range("L126").Value = 14914.2861
$ws.Range("M126").Value = -7271.75
$ws.Range("N126").Value = -19854.2861
$ws.Range("H132").Value = 3960.0967
$ws.Range("I132").Value = 3687.8462
$ws.Range("K132").Value = 11063.5386
$ws.Range("M132").Value = -8533.5386

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2573.4375
$ws.Range("I61").Value = 2508.1
$ws.Range("J61").Value = 2682.3333
$ws.Range("K61").Value = 2508.1
$ws.Range("L61").Value = 2682.3333
$ws.Range("M61").Value = -2306.1
$ws.Range("N61").Value = -3086.3333
$ws.Range("H113").Value = 2573.4375
$ws.Range("I113").Value = 2508.1
$ws.Range("J113").Value = 2682.3333
$ws.Range("K113").Value = 2508.1
$ws.Range("L113").Value = 2682.3333
$ws.Range("M113").Value = -338.0999999999999
$ws.Range("N113").Value = -7022.3333
$ws.Range("H136").Value = 81874
$ws.Range("I136").Value = 2262.6667
$ws.Range("K136").Value = 6788.000100000001
$ws.Range("M136").Value = -4238.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 15150.5
$ws.Range("J82").Value = 15150.5
$ws.Range("L82").Value = 15150.5
$ws.Range("N82").Value = -15916.5
$ws.Range("H85").Value = 15150.5
$ws.Range("J85").Value = 15150.5
$ws.Range("L85").Value = 15150.5
$ws.Range("N85").Value = -17802.5
$ws.Range("H132").Value = 3671
$ws.Range("I132").Value = 2400.1794
$ws.Range("K132").Value = 7200.5382
$ws.Range("M132").Value = -4670.5382
